$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# --- Players sheet: rows updated for 2026-01-28 final box scores ---
# Row 3
$ws1.Range("C3").Value = 'No'
$ws1.Range("D3").Value = 'Devin McGlockton'
$ws1.Range("E3").Value = 'VAN'
$ws1.Range("F3").Value = 'UK@VAN'
$ws1.Range("H3").Value = 20
$ws1.Range("I3").Value = 11
$ws1.Range("J3").Value = 12
$ws1.Range("K3").Value = 2
$ws1.Range("M3").Value = 1
$ws1.Range("N3").Value = 3
$ws1.Range("O3").Value = 32
# Row 4
$ws1.Range("D4").Value = 'Derrion Reid'
$ws1.Range("E4").Value = 'OU'
$ws1.Range("F4").Value = 'ARK@OU'
$ws1.Range("H4").Value = 12
$ws1.Range("I4").Value = 12
$ws1.Range("J4").Value = 6
$ws1.Range("K4").Value = 0
$ws1.Range("M4").Value = 0
$ws1.Range("N4").Value = 2
$ws1.Range("O4").Value = 33
# Row 7
$ws1.Range("D7").Value = 'Pablo Tamba'
$ws1.Range("E7").Value = 'LSU'
$ws1.Range("F7").Value = 'MSST@LSU'
$ws1.Range("G7").Value = 'Final'
$ws1.Range("H7").Value = 15
$ws1.Range("I7").Value = 10
$ws1.Range("J7").Value = 5
$ws1.Range("K7").Value = 2
$ws1.Range("L7").Value = 3
$ws1.Range("M7").Value = 4
$ws1.Range("N7").Value = 1
$ws1.Range("O7").Value = 35
# Row 8
$ws1.Range("D8").Value = 'Nate Ament'
$ws1.Range("E8").Value = 'TENN'
$ws1.Range("F8").Value = 'TENN@UGA'
$ws1.Range("G8").Value = 'Final/OT'
$ws1.Range("H8").Value = 13
$ws1.Range("I8").Value = 19
$ws1.Range("J8").Value = 6
$ws1.Range("K8").Value = 3
$ws1.Range("M8").Value = 0
$ws1.Range("N8").Value = 2
$ws1.Range("O8").Value = 38
# Row 9
$ws1.Range("D9").Value = 'Matas Vokietaitis'
$ws1.Range("H9").Value = 12
$ws1.Range("I9").Value = 12
$ws1.Range("L9").Value = 0
$ws1.Range("M9").Value = 1
$ws1.Range("N9").Value = 4
$ws1.Range("O9").Value = 30
# Row 10
$ws1.Range("C10").Value = 'Yes'
$ws1.Range("D10").Value = 'Tramon Mark'
$ws1.Range("E10").Value = 'TEX'
$ws1.Range("F10").Value = 'TEX@AUB'
$ws1.Range("H10").Value = -1
$ws1.Range("I10").Value = 4
$ws1.Range("K10").Value = 1
$ws1.Range("L10").Value = 2
$ws1.Range("M10").Value = 0
$ws1.Range("N10").Value = 2
$ws1.Range("O10").Value = 22
# Row 45
$ws1.Range("D45").Value = 'Malachi Moreno'
$ws1.Range("E45").Value = 'UK'
$ws1.Range("F45").Value = 'UK@VAN'
$ws1.Range("H45").Value = 11
$ws1.Range("I45").Value = 4
$ws1.Range("J45").Value = 8
$ws1.Range("K45").Value = 1
$ws1.Range("M45").Value = 1
$ws1.Range("O45").Value = 26
# Row 46
$ws1.Range("C46").Value = 'No'
$ws1.Range("D46").Value = 'Xzayvier Brown'
$ws1.Range("E46").Value = 'OU'
$ws1.Range("F46").Value = 'ARK@OU'
$ws1.Range("H46").Value = 12
$ws1.Range("I46").Value = 13
$ws1.Range("J46").Value = 6
$ws1.Range("K46").Value = 3
$ws1.Range("M46").Value = 0
$ws1.Range("O46").Value = 38
# Row 48
$ws1.Range("C48").Value = 'Yes'
$ws1.Range("D48").Value = 'Somtochukwu Cyril'
$ws1.Range("E48").Value = 'UGA'
$ws1.Range("F48").Value = 'TENN@UGA'
$ws1.Range("G48").Value = 'Final/OT'
$ws1.Range("H48").Value = 9
$ws1.Range("I48").Value = 6
$ws1.Range("J48").Value = 3
$ws1.Range("K48").Value = 0
$ws1.Range("L48").Value = 2
$ws1.Range("M48").Value = 3
$ws1.Range("N48").Value = 1
$ws1.Range("O48").Value = 32
# Row 49
$ws1.Range("C49").Value = 'Yes'
$ws1.Range("D49").Value = 'Mike Nwoko'
$ws1.Range("E49").Value = 'LSU'
$ws1.Range("F49").Value = 'MSST@LSU'
$ws1.Range("G49").Value = 'Final'
$ws1.Range("H49").Value = 3
$ws1.Range("I49").Value = 4
$ws1.Range("J49").Value = 1
$ws1.Range("K49").Value = 0
$ws1.Range("L49").Value = 1
$ws1.Range("N49").Value = 0
$ws1.Range("O49").Value = 17
# Row 50
$ws1.Range("D50").Value = 'Kevin Overton'
$ws1.Range("E50").Value = 'AUB'
$ws1.Range("F50").Value = 'TEX@AUB'
$ws1.Range("H50").Value = 25
$ws1.Range("I50").Value = 25
$ws1.Range("J50").Value = 1
$ws1.Range("K50").Value = 1
$ws1.Range("L50").Value = 1
$ws1.Range("N50").Value = 0
$ws1.Range("O50").Value = 33
# Row 51
$ws1.Range("D51").Value = 'J.P. Estrella'
$ws1.Range("E51").Value = 'TENN'
$ws1.Range("H51").Value = 20
$ws1.Range("I51").Value = 17
$ws1.Range("J51").Value = 9
$ws1.Range("K51").Value = 1
$ws1.Range("L51").Value = 0
$ws1.Range("M51").Value = 0
$ws1.Range("O51").Value = 31
# Row 52
$ws1.Range("D52").Value = 'Xaivian Lee'
$ws1.Range("E52").Value = 'FLA'
$ws1.Range("F52").Value = 'FLA@SC'
$ws1.Range("H52").Value = 14
$ws1.Range("I52").Value = 6
$ws1.Range("J52").Value = 4
$ws1.Range("K52").Value = 9
$ws1.Range("L52").Value = 2
$ws1.Range("O52").Value = 26
# Row 53
$ws1.Range("D53").Value = 'Jordan Pope'
$ws1.Range("E53").Value = 'TEX'
$ws1.Range("F53").Value = 'TEX@AUB'
$ws1.Range("H53").Value = 8
$ws1.Range("I53").Value = 12
$ws1.Range("J53").Value = 2
$ws1.Range("K53").Value = 3
$ws1.Range("L53").Value = 0
$ws1.Range("N53").Value = 2
$ws1.Range("O53").Value = 35
# Row 59
$ws1.Range("D59").Value = 'Marcus Millender'
$ws1.Range("H59").Value = 22
$ws1.Range("I59").Value = 19
$ws1.Range("J59").Value = 2
$ws1.Range("K59").Value = 4
$ws1.Range("L59").Value = 0
$ws1.Range("N59").Value = 1
$ws1.Range("O59").Value = 28
# Row 60
$ws1.Range("D60").Value = 'Kobe Knox'
$ws1.Range("E60").Value = 'SC'
$ws1.Range("F60").Value = 'FLA@SC'
$ws1.Range("H60").Value = 7
$ws1.Range("I60").Value = 6
$ws1.Range("J60").Value = 4
$ws1.Range("K60").Value = 0
$ws1.Range("L60").Value = 1
$ws1.Range("M60").Value = 1
$ws1.Range("O60").Value = 26
# Row 61
$ws1.Range("D61").Value = 'Blue Cain'
$ws1.Range("H61").Value = 9
$ws1.Range("I61").Value = 9
$ws1.Range("J61").Value = 4
$ws1.Range("K61").Value = 3
$ws1.Range("L61").Value = 1
$ws1.Range("N61").Value = 2
$ws1.Range("O61").Value = 35
# Row 62
$ws1.Range("D62").Value = 'Josh Hubbard'
$ws1.Range("E62").Value = 'MSST'
$ws1.Range("F62").Value = 'MSST@LSU'
$ws1.Range("H62").Value = 9
$ws1.Range("I62").Value = 15
$ws1.Range("J62").Value = 2
$ws1.Range("K62").Value = 3
$ws1.Range("L62").Value = 0
$ws1.Range("M62").Value = 0
$ws1.Range("O62").Value = 32
# Row 66
$ws1.Range("C66").Value = 'Yes'
$ws1.Range("D66").Value = 'Jalen Washington'
$ws1.Range("E66").Value = 'VAN'
$ws1.Range("F66").Value = 'UK@VAN'
$ws1.Range("H66").Value = 4
$ws1.Range("I66").Value = 0
$ws1.Range("J66").Value = 9
$ws1.Range("L66").Value = 0
$ws1.Range("N66").Value = 1
$ws1.Range("O66").Value = 17
# Row 67
$ws1.Range("D67").Value = 'Karter Knox'
$ws1.Range("E67").Value = 'ARK'
$ws1.Range("F67").Value = 'ARK@OU'
$ws1.Range("H67").Value = 12
$ws1.Range("I67").Value = 11
$ws1.Range("J67").Value = 2
$ws1.Range("L67").Value = 1
$ws1.Range("N67").Value = 0
$ws1.Range("O67").Value = 26

# --- OwnerTotals sheet: re-ranked by updated starter_pooh_total ---
# Row 2
$ws2.Range("B2").Value = 98
# Row 5
$ws2.Range("A5").Value = 'Booz'
$ws2.Range("B5").Value = 60
# Row 6
$ws2.Range("A6").Value = 'CDL'
$ws2.Range("B6").Value = 59
# Row 7
$ws2.Range("A7").Value = 'Tar'
$ws2.Range("B7").Value = 49
$ws2.Range("C7").Value = 5
# Row 8
$ws2.Range("A8").Value = 'Mark'
$ws2.Range("B8").Value = 46
